$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) - force text format to preserve exact string formatting
# (values like "1.001" or "88.10" would otherwise be auto-converted to numbers)
$priceUpdates = @{
    'D2' = '22.403.28'
    'D3' = '1.563.62'
    'D4' = '1.001'
    'D6' = '285.67'
    'D7' = '0.3641'
    'D8' = '48.37'
    'D9' = '0.3332'
    'D10' = '1.126'
    'D11' = '0.07407'
    'D12' = '1.001'
    'D13' = '20.77'
    'D14' = '5.924'
    'D15' = '6.882'
    'D16' = '1.563.42'
    'D17' = '0.00001104'
    'D18' = '88.10'
    'D21' = '6.350'
    'D22' = '16.07'
    'D23' = '11.96'
    'D24' = '22.391.92'
    'D25' = '2.414'
    'D27' = '149.76'
    'D28' = '19.41'
    'D29' = '4.999'
    'D30' = '123.12'
    'D31' = '1.737.40'
    'D32' = '1.065'
    'D33' = '6.118'
    'D34' = '1.992'
    'D35' = '9.580'
    'D36' = '0.08241'
    'D37' = '0.02386'
    'D38' = '1.305'
    'D39' = '0.06374'
    'D41' = '5.326'
    'D42' = '11.10'
    'D43' = '0.6067'
    'D44' = '1.0000'
    'D46' = '3.759'
    'D47' = '0.5739'
    'D48' = '2.007'
    'D49' = '124.60'
    'D51' = '0.07215'
}
foreach ($cellRef in $priceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$cellRef]
    $cell.Style = "Normal"
}

# Volume(1h) column (E) - plain percentage text, always non-numeric due to % and padding
$volumeUpdates = @{
    'E2' = '  -0.19%  '
    'E3' = '  -0.62%  '
    'E4' = '  -0.26%  '
    'E5' = '  -0.23%  '
    'E6' = '  -1.87%  '
    'E7' = '  -2.60%  '
    'E8' = '  -3.08%  '
    'E9' = '  -2.00%  '
    'E10' = '  -0.98%  '
    'E11' = '  -1.92%  '
    'E12' = '  -0.26%  '
    'E13' = '  -2.72%  '
    'E14' = '  -1.04%  '
    'E15' = '  -0.69%  '
    'E16' = '  -0.20%  '
    'E18' = '  -3.11%  '
    'E19' = '  -0.56%  '
    'E20' = '  -0.16%  '
    'E21' = '  +1.57%  '
    'E22' = '  -1.95%  '
    'E23' = '  -1.43%  '
    'E24' = '  -0.27%  '
    'E25' = '  +3.14%  '
    'E26' = '  -1.43%  '
    'E27' = '  +0.68%  '
    'E28' = '  -3.57%  '
    'E30' = '  -2.18%  '
    'E31' = '  -0.41%  '
    'E32' = '  +1.35%  '
    'E33' = '  +0.07%  '
    'E34' = '  +0.68%  '
    'E35' = '  -2.41%  '
    'E36' = '  -2.07%  '
    'E37' = '  -3.11%  '
    'E38' = '  -5.09%  '
    'E39' = '  -2.40%  '
    'E40' = '  -3.84%  '
    'E41' = '  -2.28%  '
    'E42' = '  -1.79%  '
    'E43' = '  -3.00%  '
    'E44' = '  +0.07%  '
    'E45' = '  -1.59%  '
    'E46' = '  -1.39%  '
    'E47' = '  -1.39%  '
    'E48' = '  -3.78%  '
    'E49' = '  -3.52%  '
    'E50' = '  -0.67%  '
    'E51' = '  -1.50%  '
}
foreach ($cellRef in $volumeUpdates.Keys) {
    $ws.Range($cellRef).Value = $volumeUpdates[$cellRef]
}
